$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 1382.65
$ws.Cells.Item(19, 9).Value = 1878.7693
$ws.Cells.Item(19, 10).Value = 461.2857
$ws.Cells.Item(19, 11).Value = 1878.7693
$ws.Cells.Item(19, 12).Value = 461.2857
$ws.Cells.Item(19, 13).Value = -1703.7693
$ws.Cells.Item(19, 14).Value = -811.2857

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 8126344
$ws.Cells.Item(33, 9).Value = 10626099
$ws.Cells.Item(33, 11).Value = 10626099
$ws.Cells.Item(33, 13).Value = -10625870

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(51, 8).Value = 6399.933
$ws.Cells.Item(51, 10).Value = 6769.154
$ws.Cells.Item(51, 12).Value = 6769.154
$ws.Cells.Item(51, 14).Value = -7737.154

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(131, 8).Value = 3943.4443
$ws.Cells.Item(131, 9).Value = 1098.8
$ws.Cells.Item(131, 11).Value = 3296.4
$ws.Cells.Item(131, 13).Value = 1743.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(134, 8).Value = 149999
$ws.Cells.Item(134, 10).Value = 149999
$ws.Cells.Item(134, 12).Value = 149999
$ws.Cells.Item(134, 14).Value = -160139

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 14943.25
$ws.Cells.Item(137, 9).Value = 6017.5
$ws.Cells.Item(137, 10).Value = 23869
$ws.Cells.Item(137, 11).Value = 18052.5
$ws.Cells.Item(137, 12).Value = 71607
$ws.Cells.Item(137, 13).Value = -15502.5
$ws.Cells.Item(137, 14).Value = -76707

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 2346.625
$ws.Cells.Item(138, 10).Value = 3336.2
$ws.Cells.Item(138, 12).Value = 10008.6
$ws.Cells.Item(138, 14).Value = -20288.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(35, 8).Value = 16734.75
$ws.Cells.Item(35, 10).Value = 26247.75
$ws.Cells.Item(35, 12).Value = 26247.75
$ws.Cells.Item(35, 14).Value = -27059.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 5898.6772
$ws.Cells.Item(45, 9).Value = 6308.593
$ws.Cells.Item(45, 10).Value = 3131.75
$ws.Cells.Item(45, 11).Value = 6308.593
$ws.Cells.Item(45, 12).Value = 3131.75
$ws.Cells.Item(45, 13).Value = -5931.593
$ws.Cells.Item(45, 14).Value = -3885.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 6619.5186
$ws.Cells.Item(61, 9).Value = 6251.9
$ws.Cells.Item(61, 10).Value = 6835.7646
$ws.Cells.Item(61, 11).Value = 6251.9
$ws.Cells.Item(61, 12).Value = 6835.7646
$ws.Cells.Item(61, 13).Value = -6039.9
$ws.Cells.Item(61, 14).Value = -7259.7646

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 16629.088
$ws.Cells.Item(74, 9).Value = 19444.277
$ws.Cells.Item(74, 10).Value = 6494.4
$ws.Cells.Item(74, 11).Value = 19444.277
$ws.Cells.Item(74, 12).Value = 6494.4
$ws.Cells.Item(74, 13).Value = -18570.277
$ws.Cells.Item(74, 14).Value = -8242.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 16629.088
$ws.Cells.Item(77, 9).Value = 19444.277
$ws.Cells.Item(77, 10).Value = 6494.4
$ws.Cells.Item(77, 11).Value = 97221.38499999999
$ws.Cells.Item(77, 12).Value = 32472
$ws.Cells.Item(77, 13).Value = -92853.38499999999
$ws.Cells.Item(77, 14).Value = -41208

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 1634.0526
$ws.Cells.Item(102, 9).Value = 1765.1428
$ws.Cells.Item(102, 11).Value = 1765.1428
$ws.Cells.Item(102, 13).Value = -143.1428000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 1849.3334
$ws.Cells.Item(122, 9).Value = 1849.375
$ws.Cells.Item(122, 10).Value = 1849.25
$ws.Cells.Item(122, 11).Value = 5548.125
$ws.Cells.Item(122, 12).Value = 5547.75
$ws.Cells.Item(122, 13).Value = -3098.125
$ws.Cells.Item(122, 14).Value = -10447.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 1056.7407
$ws.Cells.Item(132, 9).Value = 1017.5
$ws.Cells.Item(132, 11).Value = 3052.5
$ws.Cells.Item(132, 13).Value = -522.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 6619.5186
$ws.Cells.Item(136, 9).Value = 6251.9
$ws.Cells.Item(136, 10).Value = 6835.7646
$ws.Cells.Item(136, 11).Value = 18755.7
$ws.Cells.Item(136, 12).Value = 20507.2938
$ws.Cells.Item(136, 13).Value = -16205.7
$ws.Cells.Item(136, 14).Value = -25607.2938

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(19, 8).Value = 0
$ws.Cells.Item(19, 9).Value = 0
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 11).Value = 0
$ws.Cells.Item(19, 12).Value = 0
$ws.Cells.Item(19, 13).ClearContents() | Out-Null
$ws.Cells.Item(19, 14).ClearContents() | Out-Null

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(54, 8).Value = 5190.8
$ws.Cells.Item(54, 9).Value = 5190.8
$ws.Cells.Item(54, 11).Value = 5190.8
$ws.Cells.Item(54, 13).Value = -4706.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 3312.1538
$ws.Cells.Item(107, 9).Value = 2170
$ws.Cells.Item(107, 10).Value = 5139.6
$ws.Cells.Item(107, 11).Value = 2170
$ws.Cells.Item(107, 12).Value = 5139.6
$ws.Cells.Item(107, 13).Value = -250
$ws.Cells.Item(107, 14).Value = -8979.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 6542.865
$ws.Cells.Item(134, 9).Value = 2624.389
$ws.Cells.Item(134, 11).Value = 7873.167
$ws.Cells.Item(134, 13).Value = -5338.167

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(42, 8).Value = 9866.4
$ws.Cells.Item(42, 9).Value = 4639
$ws.Cells.Item(42, 10).Value = 13351.333
$ws.Cells.Item(42, 11).Value = 4639
$ws.Cells.Item(42, 12).Value = 13351.333
$ws.Cells.Item(42, 13).Value = -4046
$ws.Cells.Item(42, 14).Value = -14537.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 3268.0227
$ws.Cells.Item(58, 9).Value = 1636.6471
$ws.Cells.Item(58, 11).Value = 1636.6471
$ws.Cells.Item(58, 13).Value = -1433.6471

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 9461.5
$ws.Cells.Item(99, 10).Value = 11052.725
$ws.Cells.Item(99, 12).Value = 11052.725
$ws.Cells.Item(99, 14).Value = -14048.725

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 9461.5
$ws.Cells.Item(126, 10).Value = 11052.725
$ws.Cells.Item(126, 12).Value = 33158.175
$ws.Cells.Item(126, 14).Value = -38098.175

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 3268.0227
$ws.Cells.Item(136, 9).Value = 1636.6471
$ws.Cells.Item(136, 11).Value = 4909.9413
$ws.Cells.Item(136, 13).Value = -2359.9413

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 62569730
$ws.Cells.Item(7, 9).Value = 125000180
$ws.Cells.Item(7, 11).Value = 375000540
$ws.Cells.Item(7, 13).Value = -375000428

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(13, 8).Value = 200
$ws.Cells.Item(13, 9).Value = 0
$ws.Cells.Item(13, 10).Value = 200
$ws.Cells.Item(13, 11).Value = 0
$ws.Cells.Item(13, 12).Value = 600
$ws.Cells.Item(13, 13).ClearContents() | Out-Null
$ws.Cells.Item(13, 14).Value = -936

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(47, 8).Value = 31
$ws.Cells.Item(47, 9).Value = 31
$ws.Cells.Item(47, 11).Value = 93
$ws.Cells.Item(47, 13).Value = 338

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(132, 8).Value = 4349672.5
$ws.Cells.Item(132, 9).Value = 1668.1428
$ws.Cells.Item(132, 10).Value = 6251924.5
$ws.Cells.Item(132, 11).Value = 15013.2852
$ws.Cells.Item(132, 12).Value = 56267320.5
$ws.Cells.Item(132, 13).Value = -12483.2852
$ws.Cells.Item(132, 14).Value = -56272380.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(133, 8).Value = 5000
$ws.Cells.Item(133, 10).Value = 0
$ws.Cells.Item(133, 12).Value = 0
$ws.Cells.Item(133, 14).ClearContents() | Out-Null

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(99, 8).Value = 35323.668
$ws.Cells.Item(99, 10).Value = 49990
$ws.Cells.Item(99, 12).Value = 49990
$ws.Cells.Item(99, 14).Value = -54482

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 5379.8
$ws.Cells.Item(122, 9).Value = 1499.5
$ws.Cells.Item(122, 10).Value = 7966.6665
$ws.Cells.Item(122, 11).Value = 4498.5
$ws.Cells.Item(122, 12).Value = 23899.9995
$ws.Cells.Item(122, 13).Value = -2048.5
$ws.Cells.Item(122, 14).Value = -28799.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 12055.047
$ws.Cells.Item(132, 9).Value = 10362.419
$ws.Cells.Item(132, 11).Value = 31087.257
$ws.Cells.Item(132, 13).Value = -28557.257

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(56, 8).Value = 26704.334
$ws.Cells.Item(56, 10).Value = 26704.334
$ws.Cells.Item(56, 12).Value = 26704.334
$ws.Cells.Item(56, 14).Value = -28086.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 2296.1614
$ws.Cells.Item(93, 9).Value = 2841.9524
$ws.Cells.Item(93, 11).Value = 2841.9524
$ws.Cells.Item(93, 13).Value = -1593.9524

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 7073.4595
$ws.Cells.Item(132, 9).Value = 6637.6333
$ws.Cells.Item(132, 11).Value = 19912.8999
$ws.Cells.Item(132, 13).Value = -17382.8999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(133, 8).Value = 99998.336
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 10).Value = 99998.336
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 12).Value = 99998.336
$ws.Cells.Item(133, 13).ClearContents() | Out-Null
$ws.Cells.Item(133, 14).Value = -105058.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 5132.025
$ws.Cells.Item(136, 9).Value = 4290.2573
$ws.Cells.Item(136, 10).Value = 11024.4
$ws.Cells.Item(136, 11).Value = 12870.7719
$ws.Cells.Item(136, 12).Value = 33073.2
$ws.Cells.Item(136, 13).Value = -10320.7719
$ws.Cells.Item(136, 14).Value = -38173.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(42, 8).Value = 72524
$ws.Cells.Item(42, 9).Value = 49999
$ws.Cells.Item(42, 11).Value = 49999
$ws.Cells.Item(42, 13).Value = -49621

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 1277.7778
$ws.Cells.Item(96, 9).Value = 882.8095
$ws.Cells.Item(96, 11).Value = 882.8095
$ws.Cells.Item(96, 13).Value = 490.1905

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 5799.037
$ws.Cells.Item(122, 9).Value = 3759.5908
$ws.Cells.Item(122, 11).Value = 11278.7724
$ws.Cells.Item(122, 13).Value = -8828.7724

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 17775.584
$ws.Cells.Item(132, 9).Value = 11344.424
$ws.Cells.Item(132, 11).Value = 34033.272
$ws.Cells.Item(132, 13).Value = -31503.272
